$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns (in this worksheet) that receive updated values in this edit,
# in order: B, C, D, E, G, J, L, M, O
$cols = @(2, 3, 4, 5, 7, 10, 12, 13, 15)

# New values for rows 2-25 (one row per data record), one value per column above
$data = @(
    @(1.776994643007185, 0.1408982152548504, 0.4918735383361792, 0.1564770308378582, 0.002550075684296477, 0.06522862455345191, 0.405987759196492, 0.4387683941067948, 7.179538403729794),
    @(1.689616518605021, 0.1291951117967471, 0.4925051583194602, 0.157553810353475, 0.002553809453376121, 0.06454635528047348, 0.4027974556828937, 0.4245932815990798, 7.214462352219215),
    @(1.636559501023214, 0.1219643690044876, 0.4930935991483096, 0.1582578767860392, 0.002556224857281407, 0.06412367458509038, 0.4009887467442113, 0.4160515212538769, 7.239769016823288),
    @(1.615088780636427, 0.119006606811098, 0.4933838972863356, 0.158555596970535, 0.002557240143025683, 0.06395048940904147, 0.4002895381556613, 0.4126115942789568, 7.251052432392555),
    @(1.611532706173875, 0.1185148017536335, 0.49343515300761, 0.1586056863603584, 0.002557410604863942, 0.06392167562488282, 0.400175724370385, 0.4120428744020685, 7.252984645315223),
    @(1.636269328627009, 0.1219245246537071, 0.49309730964994, 0.1582618481659646, 0.002556238424077702, 0.0641213427429399, 0.4009791635675626, 0.4160049632605336, 7.239917259350932),
    @(1.746744171047396, 0.1368724011239806, 0.492049713645244, 0.1568394055259894, 0.002551337647578644, 0.06499416358806442, 0.4048566358249133, 0.4338473561587364, 7.190777943422319),
    @(1.968054232176712, 0.1658235747168817, 0.4915852332092925, 0.1543898602950398, 0.002542697633098297, 0.0666756535079891, 0.4136481927077398, 0.4701127769596667, 7.12510489631191),
    @(2.133462596131039, 0.1868695333655808, 0.4922111265890976, 0.1527964017469925, 0.002536935227196767, 0.06789248772002665, 0.4208279247550308, 0.4975285936534846, 7.095620701113916),
    @(2.209315824143232, 0.1963944522834424, 0.4927054821025649, 0.1521160703977769, 0.002534439563470388, 0.06844199450290489, 0.4242499742841517, 0.5101671622023005, 7.086295245658619),
    @(2.238126040984071, 0.1999941453717895, 0.4929227847963773, 0.1518648365223534, 0.002533512494030058, 0.06864949266011067, 0.4255681645351359, 0.5149769084059699, 7.083352607983727),
    @(2.231917428482802, 0.1992192089338118, 0.4928746467195282, 0.151918660100745, 0.002533711356549183, 0.06860483049610266, 0.4252832767976003, 0.5139399886459231, 7.083960158700563),
    @(2.211684336361941, 0.196690746041952, 0.4927227567075505, 0.152095273182522, 0.002534362932807039, 0.0684590773316387, 0.424357975495937, 0.5105623872668446, 7.086041348739229),
    @(2.199302179530207, 0.1951410493887806, 0.4926336384278613, 0.1522042859741646, 0.002534764381863985, 0.06836972249428896, 0.4237941076268186, 0.508496601283575, 7.087392834792468),
    @(2.128517543238729, 0.1862460594640538, 0.4921830338790727, 0.1528417582424479, 0.002537100846174913, 0.06785649426867835, 0.4206074160996138, 0.4967059748961873, 7.096312468610336),
    @(2.085248357556395, 0.1807766212113791, 0.4919602632477194, 0.1532442261999183, 0.002538566318617492, 0.06754060590505517, 0.4186923619709972, 0.4895154191705799, 7.102831875182488),
    @(2.060418377615861, 0.1776261418362992, 0.4918518657716646, 0.1534799082445284, 0.002539421054813179, 0.06735853571761297, 0.4176055579624034, 0.485395331017159, 7.106966298227263),
    @(2.052021253255475, 0.1765586574171323, 0.4918185549844338, 0.1535604266777595, 0.002539712489215733, 0.06729682492486688, 0.4172401096800229, 0.4840030481287911, 7.108432173626341),
    @(2.089848513511754, 0.1813593296949989, 0.4919819354571189, 0.1532009489143995, 0.002538409092270988, 0.06757427209853972, 0.4188947036311248, 0.490279239754571, 7.102098057935677),
    @(2.217624954797543, 0.1974336130550398, 0.4927665539029533, 0.1520432242352401, 0.002534171061864787, 0.0685019045907751, 0.4246291534936404, 0.5115538265972717, 7.085414067391355),
    @(2.301636405065324, 0.2078971311347289, 0.4934547789396788, 0.1513238410630571, 0.00253150604966968, 0.06910473451950594, 0.4285070880636965, 0.525596593731251, 7.077942032844703),
    @(2.256752358419362, 0.202316439738695, 0.4930714227110826, 0.151704384367723, 0.002532918858078399, 0.06878330929803411, 0.4264254835227206, 0.5180890939020273, 7.081615647151239),
    @(2.087768639467527, 0.1810959059529011, 0.4919720761582624, 0.1532205011867571, 0.002538480136257157, 0.06755905305020704, 0.4188031807723576, 0.4899338732857785, 7.102428613350867),
    @(1.907687770157338, 0.1580306514730694, 0.4915408986418299, 0.1550162429698378, 0.002544931738464153, 0.06622401065471095, 0.4111430497333828, 0.4601660512374934, 7.139579855269602)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $values[$j]
    }
}

$wb.Save()
